$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 11 (G + H change to 1)
foreach ($r in @(3, 11)) {
    $ws.Cells.Item($r, 7).Value = 1   # column G
    $ws.Cells.Item($r, 8).Value = 1   # column H
}

# Rows 4, 5, 12 (D + E change to 1)
foreach ($r in @(4, 5, 12)) {
    $ws.Cells.Item($r, 4).Value = 1   # column D
    $ws.Cells.Item($r, 5).Value = 1   # column E
}

# Rows 6-10, 13-18 (H changes to 1)
foreach ($r in @(6, 7, 8, 9, 10, 13, 14, 15, 16, 17, 18)) {
    $ws.Cells.Item($r, 8).Value = 1   # column H
}
